$d = $word.ActiveDocument

function FindRange($text) {
    $r = $d.Content
    $found = $r.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND:" $text
    }
    return $r
}

# ============ EDIT 1: Insert "Round-off errors" (Heading3) before "Letters and words" ============
$r1 = FindRange("Letters and words")
$r1.Collapse(1)
$r1.InsertBefore("Round-off errors`r")
$newPara1 = $r1.Paragraphs(1).Previous()
$newPara1.Style = "Heading 3"

Write-Host "Edit1 done:" $newPara1.Range.Text

# ============ EDIT 2: Build "Java Methods" ... "Java Classes" section ============
# The bookmark "_GoBack" currently sits alone in an empty Heading1 paragraph.
$bm = $d.Bookmarks("_GoBack")

# -- Insert "Java Methods" text right before the bookmark (same paragraph) --
$p1 = $bm.Range.Duplicate
$p1.Collapse(1)
$p1.InsertBefore("Java Methods")

# -- Insert 2 new Heading2 paragraphs before the (still-empty) bookmark paragraph --
$bm = $d.Bookmarks("_GoBack")
$p2 = $bm.Range.Duplicate
$p2.Collapse(1)
$p2.InsertBefore("What is a method?`rSequential execution`r")

# -- Insert " first method!" right after the bookmark, then "Your" right before it --
$bm = $d.Bookmarks("_GoBack")
$p3 = $bm.Range.Duplicate
$p3.Collapse(0)
$p3.InsertAfter(" first method!")

$bm = $d.Bookmarks("_GoBack")
$p4 = $bm.Range.Duplicate
$p4.Collapse(1)
$p4.InsertBefore("Your")

# Fix styles of the newly created paragraphs (they all inherited Heading1 so far,
# which is already correct for "Java Methods"; the rest need Heading2).
$bm = $d.Bookmarks("_GoBack")
$bmParaRange = $bm.Range.Duplicate
$bmParaRange.Expand(4)
$bmPara = $bmParaRange.Paragraphs(1)
$bmPara.Style = "Heading 2"                     # "Your first method!"
$bmPara.Previous().Style = "Heading 2"          # "Sequential execution"
$bmPara.Previous().Previous().Style = "Heading 2"   # "What is a method?"

# -- Insert the rest of the new paragraphs right after the bookmark's paragraph --
$bm = $d.Bookmarks("_GoBack")
$bmParaRange = $bm.Range.Duplicate
$bmParaRange.Expand(4)
$nextParaStart = $bmParaRange.End
$p5 = $d.Range($nextParaStart, $nextParaStart)
$p5.InsertBefore("Parameters`rFormal parameters and Actual parameters`rJava Classes`rWhat is a class?`rYour first class!`r`r")

# Re-apply explicit styles (also clears inherited numPr/ind overrides picked up
# from the paragraph that followed the insertion point).
$afterBmPara = $bmPara.Next()
$afterBmPara.Style = "Heading 2"              # Parameters
$p = $afterBmPara.Next()
$p.Style = "Heading 3"                        # Formal parameters and Actual parameters
$p = $p.Next()
$p.Style = "Heading 1"                        # Java Classes
$p = $p.Next()
$p.Style = "Heading 2"                        # What is a class?
$p = $p.Next()
$p.Style = "Heading 2"                        # Your first class!
$p = $p.Next()
$p.Style = "Heading 1"                        # (new empty paragraph)

Write-Host "Edit2 done"
